$d = $word.ActiveDocument

# NOTE: paragraph indices below refer to the ORIGINAL (pre-edit) document.
# We process edits from the bottom of the document upward so that each
# edit's index references stay valid (earlier / lower-numbered paragraphs
# are untouched by edits made further down in the document).

# ---------------------------------------------------------------------
# 1) "HIVE TEAM: DEVELOPMENT" member block - remove the four paragraphs
#    (emelia / auditor / Nitego / Senior QA Tester) that followed
#    "Release Coordinator", leaving "Jazz" / "Release Coordinator" /
#    "clot" / "security" ... contiguous.
# ---------------------------------------------------------------------
$pRelStart = $d.Paragraphs.Item(50)
$pRelEnd   = $d.Paragraphs.Item(53)
$rngRel = $d.Range($pRelStart.Range.Start, $pRelEnd.Range.End)
$rngRel.Delete()

# ---------------------------------------------------------------------
# 2) "HIVE TEAM: OUTREACH" member block - remove the four paragraphs
#    (Outreach Support / Carlos Santiago / Outreach Support / Emilio)
#    that followed "Semptly", leaving "Semptly" / "Outreach Support" /
#    "misachasu" ... contiguous.
# ---------------------------------------------------------------------
$pSemStart = $d.Paragraphs.Item(16)
$pSemEnd   = $d.Paragraphs.Item(19)
$rngSem = $d.Range($pSemStart.Range.Start, $pSemEnd.Range.End)
$rngSem.Delete()

# ---------------------------------------------------------------------
# 3) Replace the "Mark Hakkarinen" member entry (4 paragraphs: name,
#    title, "Editor of ..." line with hyperlink, "Email" hyperlink)
#    with the new translated entry for "LilyDaVine" (2 paragraphs: name
#    heading, "Outreach Support" title).
# ---------------------------------------------------------------------
$pNameStart = $d.Paragraphs.Item(11)
$pNameEnd   = $d.Paragraphs.Item(14)
$rngName = $d.Range($pNameStart.Range.Start, $pNameEnd.Range.End)

$newEntryXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading3"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="105" w:beforeAutospacing="0" w:after="120" w:afterAutospacing="0" w:line="264" w:lineRule="atLeast"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="Exo" w:hAnsi="Exo"/><w:b w:val="0"/><w:bCs w:val="0"/><w:color w:val="343434"/><w:sz w:val="37"/><w:szCs w:val="37"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Exo" w:hAnsi="Exo"/><w:b w:val="0"/><w:bCs w:val="0"/><w:color w:val="343434"/><w:sz w:val="37"/><w:szCs w:val="37"/></w:rPr><w:t>LilyDaVine</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans" w:eastAsia="Times New Roman"/><w:color w:val="3B3B3B"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">Outreach Support</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$rngName.InsertXML($newEntryXml)
